$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that have an "Additional Effort [h]" value in column C that needs to be
# folded into column B (Effort [h]) before the column is removed.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $add = $ws.Cells.Item($r, 3).Value()
    if ($add -ne $null) {
        $base = $ws.Cells.Item($r, 2).Value()
        if ($base -eq $null) { $base = 0 }
        $ws.Cells.Item($r, 2).Value = $base + $add
    }
}

# Remove the now-merged "Additional Effort [h]" column entirely (this also
# shifts the comment column D into C).
$ws.Columns.Item(3).Delete()
